$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that would otherwise be parsed as numbers by Excel,
# since the source data keeps them as literal text strings (e.g. "1.001", "4.500").
foreach ($row in @(4,5,6,7,8,9,10,12,13,14,15,16,17,18,20,21,23,24,25,26,29,30,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,47,48,49,50,51)) {
    $ws.Range("D$row").NumberFormat = "@"
}

$ws.Range("D2").Value = "25.838.60"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.736.28"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "230.11"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5180"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "0.2749"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").Value = "39.32"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "0.06122"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "1.733.39"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "0.07058"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "15.14"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "0.6352"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "4.500"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "76.87"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "25.833.65"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "11.45"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "0.000006634"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "1.955.51"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Value = "4.130"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "8.715"
$ws.Range("D25").Value = "5.128"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "139.57"
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D29").Value = "1.776"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "101.86"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "0.08302"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "3.685"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "3.476"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").Value = "0.04491"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "2.616"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "0.9730"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "0.6133"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "2.654"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").Value = "1.941"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "100.26"
$ws.Range("D43").Value = "0.3815"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D46").Value = "0.05376"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "0.1122"
$ws.Range("D48").Value = "6.218"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").Value = "52.88"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "29.91"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "7.599"
$ws.Range("E51").Value = "  +2.21%  "

# Rows 44/45: FraxShare and TrustWalletToken swap ranking positions with updated data
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "4.999"
$ws.Range("E44").Value = "  +2.02%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.7216"
$ws.Range("E45").Value = "  -3.80%  "
